# Insert a new weekly price-report row at row 44 (shifting the existing
# rows 44-82 down to 45-83, preserving all their data unchanged) and
# populate the newly inserted row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(44).Insert()

$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44586
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = 100112026
$ws.Range("G44").Value = "Haba"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 80
$ws.Range("K44").Value = 22000
$ws.Range("L44").Value = 22000
$ws.Range("M44").Value = 22000
$ws.Range("N44").Value = "`$/saco 25 kilos"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 880
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
